$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the simulated "time"/"Energy" series (columns B and C) for rows 2-49
# with results from the improved (priority-queue based) computation.
$ws.Range("B2").Value = 0.9444217277969167
$ws.Range("C2").Value = 0.7596385166932943
$ws.Range("B3").Value = 2.774620064875066
$ws.Range("C3").Value = 1.567219113513667
$ws.Range("B4").Value = 2.842811268035476
$ws.Range("C4").Value = 2.48282065084712
$ws.Range("B5").Value = 6.317145538704567
$ws.Range("C5").Value = 3.264369499221528
$ws.Range("B6").Value = 14.28682616887567
$ws.Range("C6").Value = 4.032842906585683
$ws.Range("B7").Value = 14.70307976002099
$ws.Range("C7").Value = 4.857131403670579
$ws.Range("B8").Value = 24.33633155721219
$ws.Range("C8").Value = 5.830088419933242
$ws.Range("B9").Value = 24.35540936530362
$ws.Range("C9").Value = 6.634145107410721
$ws.Range("B10").Value = 24.55520156492015
$ws.Range("C10").Value = 7.532931176900971
$ws.Range("B11").Value = 25.58737130503432
$ws.Range("C11").Value = 8.617656791803665
$ws.Range("B12").Value = 27.56505196871588
$ws.Range("C12").Value = 9.775384836645197
$ws.Range("B13").Value = 27.59744921442942
$ws.Range("C13").Value = 10.60364008882838
$ws.Range("B14").Value = 29.19384910424169
$ws.Range("C14").Value = 11.51125620987772
$ws.Range("B15").Value = 29.95059616164134
$ws.Range("C15").Value = 12.31078490084895
$ws.Range("B16").Value = 31.11240284493924
$ws.Range("C16").Value = 13.17359041439532
$ws.Range("B17").Value = 39.32139463291346
$ws.Range("C17").Value = 14.08021841191362
$ws.Range("B18").Value = 41.06528662456996
$ws.Range("C18").Value = 15.11857542439947
$ws.Range("B19").Value = 41.12681888971752
$ws.Range("C19").Value = 15.94254763311054
$ws.Range("B20").Value = 43.9714449547551
$ws.Range("C20").Value = 16.82672726633413
$ws.Range("B21").Value = 44.13687137843229
$ws.Range("C21").Value = 17.76365235960627
$ws.Range("B22").Value = 44.9506947026916
$ws.Range("C22").Value = 18.54985545842045
$ws.Range("B23").Value = 50.92332393873234
$ws.Range("C23").Value = 19.47020674008733
$ws.Range("B24").Value = 51.36930311469564
$ws.Range("C24").Value = 20.20014150814491
$ws.Range("B25").Value = 52.51444853689464
$ws.Range("C25").Value = 21.15309291617465
$ws.Range("B26").Value = 52.61457263523537
$ws.Range("C26").Value = 22.1148503524597
$ws.Range("B27").Value = 54.19806292821032
$ws.Range("C27").Value = 23.32560865525048
$ws.Range("B28").Value = 56.4623842614324
$ws.Range("C28").Value = 24.33288942349829
$ws.Range("B29").Value = 56.5093199449222
$ws.Range("C29").Value = 25.13269264296799
$ws.Range("B30").Value = 57.471175849984
$ws.Range("C30").Value = 25.95504860793136
$ws.Range("B31").Value = 60.37867715277268
$ws.Range("C31").Value = 26.75094353749739
$ws.Range("B32").Value = 60.44543289110976
$ws.Range("C32").Value = 27.52147650961148
$ws.Range("B33").Value = 68.3327834101155
$ws.Range("C33").Value = 28.40874845515717
$ws.Range("B34").Value = 68.44928422432412
$ws.Range("C34").Value = 29.4418428747889
$ws.Range("B35").Value = 69.95736642941638
$ws.Range("C35").Value = 30.23420146340424
$ws.Range("B36").Value = 70.35226861512035
$ws.Range("C36").Value = 31.12706738342384
$ws.Range("B37").Value = 73.08823877089569
$ws.Range("C37").Value = 31.90897906612092
$ws.Range("B38").Value = 78.23280621695574
$ws.Range("C38").Value = 32.78346059655551
$ws.Range("B39").Value = 78.30158476755604
$ws.Range("C39").Value = 33.62149765349496
$ws.Range("B40").Value = 79.38493040228788
$ws.Range("C40").Value = 34.39942600929621
$ws.Range("B41").Value = 80.11246476280844
$ws.Range("C41").Value = 35.21398455304369
$ws.Range("B42").Value = 80.14642671395357
$ws.Range("C42").Value = 36.05970035090213
$ws.Range("B43").Value = 80.20121364511719
$ws.Range("C43").Value = 37.01599325989163
$ws.Range("B44").Value = 80.53196129647057
$ws.Range("C44").Value = 37.80106765351571
$ws.Range("B45").Value = 80.60161889964485
$ws.Range("C45").Value = 38.76069991282043
$ws.Range("B46").Value = 88.09664982805778
$ws.Range("C46").Value = 39.93725002330507
$ws.Range("B47").Value = 90.70728185263094
$ws.Range("C47").Value = 40.85446698094874
$ws.Range("B48").Value = 95.91199811428393
$ws.Range("C48").Value = 41.58407775096931
$ws.Range("B49").Value = 95.97101012976707
$ws.Range("C49").Value = 42.50143184613968

# Add new row 50
$ws.Range("A49").Copy()
$ws.Range("A50").PasteSpecial(-4122)
$ws.Range("A50").Value = 48
$ws.Range("B50").Value = 98.78972240088119
$ws.Range("C50").Value = 43.26056095065733

"done"